$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1103.7
$ws.Range("C2").Value = 125

$ws.Range("B3").Value = 541
$ws.Range("C3").Value = 125

$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 125

$ws.Range("B5").Value = 219
$ws.Range("C5").Value = 125

$ws.Range("B6").Value = 92
$ws.Range("C6").Value = 125

$ws.Range("B7").Value = 196
$ws.Range("C7").Value = 125

$ws.Range("B10").Value = 915
$ws.Range("C10").Value = 125

$ws.Range("B11").Value = 372
$ws.Range("C11").Value = 125

$ws.Range("B12").Value = 1019.2
$ws.Range("C12").Value = 125

$ws.Range("B13").Value = 766
$ws.Range("C13").Value = 125

$ws.Range("B14").Value = 1157
$ws.Range("C14").Value = 125

$ws.Range("B15").Value = 309
$ws.Range("C15").Value = 125

$ws.Range("B17").Value = 191
$ws.Range("C17").Value = 125

$ws.Range("B18").Value = 4
$ws.Range("C18").Value = 125

$ws.Range("B19").Value = 22
$ws.Range("C19").Value = 125

$ws.Range("B20").Value = 108
$ws.Range("C20").Value = 125
